# Updated cryptos list with latest Price / Volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells are forced to text via a leading apostrophe so that
# Excel does not reinterpret dotted "thousands-style" numbers (e.g. 289.90,
# 1.149) as floating point values; the style is then reset to Normal so no
# numeric/quote-prefix formatting gets stamped onto the cell.

$ws.Range("D2").Value = "'22.440.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "'1.570.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "'289.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'0.3693"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.46%  "

$ws.Range("D8").Value = "'49.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "'0.3382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").Value = "'1.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").Value = "'0.07558"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "'21.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").Value = "'6.025"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").Value = "'6.989"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("D16").Value = "'1.571.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "'0.00001122"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "'90.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").Value = "'6.363"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.20%  "

$ws.Range("D22").Value = "'16.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").Value = "'12.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.98%  "

$ws.Range("D24").Value = "'22.447.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("D25").Value = "'2.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").Value = "'2.669"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").Value = "'20.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "

$ws.Range("D28").Value = "'149.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("D29").Value = "'5.054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.20%  "

$ws.Range("D30").Value = "'125.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("D31").Value = "'1.747.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").Value = "'1.063"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.99%  "

$ws.Range("D33").Value = "'6.231"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.99%  "

$ws.Range("D34").Value = "'2.015"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").Value = "'9.817"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").Value = "'0.08381"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("D37").Value = "'0.02476"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").Value = "'1.350"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.88%  "

$ws.Range("D39").Value = "'0.2302"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("D40").Value = "'0.06571"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").Value = "'5.420"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").Value = "'11.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.39%  "

$ws.Range("D43").Value = "'0.6253"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("D44").Value = "'14.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.39%  "

$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").Value = "'3.800"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").Value = "'0.5873"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.40%  "

$ws.Range("D48").Value = "'2.073"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("D49").Value = "'127.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.31%  "

$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").Value = "'0.07308"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.16%  "
